$d = $word.ActiveDocument

function New-WordPkgXml($bodyFrag) {
    return '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyFrag + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- 1. Strip the yellow highlight from the last ("TO DO: UPDATE REPORT...") run ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.HighlightColorIndex = 0

# --- 2. Open up a fresh (empty) trailing paragraph to host the new content ---
$lastPara.Range.InsertParagraphAfter()

# --- Paragraph: bold date heading "17/11/2023" ---
$paraDate = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>17/11/2023</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertXML((New-WordPkgXml $paraDate))

# --- Paragraph: "Wrote scripts for colocalisation UBCLung eQTL ! ..." ---
$paraBody = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Wrote scripts for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>colocalisation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>UBCLung</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>eQTL</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> !</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> 003_submit_coloc_susie_lung_eQTL.sh and 003_run_</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>coloc_susie_lung_eQTL</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.r</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>, and 004_concat_coloc_susie_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>ubclung.R</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertXML((New-WordPkgXml $paraBody))

# --- Paragraph: "TO DO: ONCE JOBS FINISHED, RUN 004_concat_coloc_susie_ubclung.R" ---
$paraTodo = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>TO DO: ONCE JOBS FINISHED, RUN 004_concat_coloc_susie_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>ubclung.R</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertXML((New-WordPkgXml $paraTodo))

# --- 3. InsertXML always leaves one trailing empty paragraph behind; drop it ---
$d.Paragraphs.Last.Range.Delete()
